$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 61999.5
$ws.Range("J95").Value = 61999.5
$ws.Range("L95").Value = 61999.5
$ws.Range("N95").Value = -67491.5
$ws.Range("H107").Value = 2998.75
$ws.Range("I107").Value = 3000
$ws.Range("J107").Value = 2998.3333
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 2998.3333
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -6838.3333
$ws.Range("H113").Value = 6408.091
$ws.Range("I113").Value = 5570
$ws.Range("K113").Value = 5570
$ws.Range("M113").Value = -2316
$ws.Range("H132").Value = 3758.7585
$ws.Range("I132").Value = 3538.3845
$ws.Range("K132").Value = 10615.1535
$ws.Range("M132").Value = -8085.1535
$ws.Range("H137").Value = 4492.577
$ws.Range("I137").Value = 1452.9131
$ws.Range("K137").Value = 4358.7393
$ws.Range("M137").Value = -1808.7393

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6184.477
$ws.Range("I32").Value = 5337.657
$ws.Range("K32").Value = 5337.657
$ws.Range("M32").Value = -5050.657
$ws.Range("H74").Value = 24392934
$ws.Range("I74").Value = 90911110
$ws.Range("K74").Value = 90911110
$ws.Range("M74").Value = -90910236
$ws.Range("H77").Value = 24392934
$ws.Range("I77").Value = 90911110
$ws.Range("K77").Value = 454555550
$ws.Range("M77").Value = -454551182
$ws.Range("H102").Value = 2867.8262
$ws.Range("I102").Value = 2358.389
$ws.Range("K102").Value = 2358.389
$ws.Range("M102").Value = -736.3890000000001
$ws.Range("H110").Value = 6282.3335
$ws.Range("I110").Value = 7082
$ws.Range("K110").Value = 7082
$ws.Range("M110").Value = -5037
$ws.Range("H122").Value = 9010828
$ws.Range("I122").Value = 891.9583
$ws.Range("K122").Value = 2675.8749
$ws.Range("M122").Value = -225.8748999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 16484.875
$ws.Range("I97").Value = 16484.875
$ws.Range("K97").Value = 16484.875
$ws.Range("M97").Value = -15493.875
$ws.Range("H99").Value = 1702
$ws.Range("I99").Value = 1702
$ws.Range("K99").Value = 1702
$ws.Range("M99").Value = -204
$ws.Range("H105").Value = 6041.5757
$ws.Range("I105").Value = 12429.7
$ws.Range("J105").Value = 3264.1304
$ws.Range("K105").Value = 12429.7
$ws.Range("L105").Value = 3264.1304
$ws.Range("M105").Value = -10682.7
$ws.Range("N105").Value = -6758.1304

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 19282.5
$ws.Range("J28").Value = 19282.5
$ws.Range("L28").Value = 19282.5
$ws.Range("N28").Value = -19772.5
$ws.Range("H37").Value = 5000
$ws.Range("J37").Value = 5000
$ws.Range("L37").Value = 5000
$ws.Range("H62").Value = 3064.7334
$ws.Range("I62").Value = 2803.2222
$ws.Range("J62").Value = 3457
$ws.Range("K62").Value = 2803.2222
$ws.Range("L62").Value = 3457
$ws.Range("M62").Value = -2179.2222
$ws.Range("N62").Value = -4705
$ws.Range("H65").Value = 3064.7334
$ws.Range("I65").Value = 2803.2222
$ws.Range("J65").Value = 3457
$ws.Range("K65").Value = 14016.111
$ws.Range("L65").Value = 17285
$ws.Range("M65").Value = -10896.111
$ws.Range("N65").Value = -23525
$ws.Range("H107").Value = 972.6429000000001
$ws.Range("I107").Value = 925.1
$ws.Range("K107").Value = 925.1
$ws.Range("M107").Value = 994.9
$ws.Range("H134").Value = 2974.5557
$ws.Range("I134").Value = 2436.1333
$ws.Range("K134").Value = 7308.3999
$ws.Range("M134").Value = -4773.3999
$ws.Range("N37").Value = -5214

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 506.45834
$ws.Range("J12").Value = 497.33334
$ws.Range("L12").Value = 1492.00002
$ws.Range("N12").Value = -1838.00002
$ws.Range("H17").Value = 348.81818
$ws.Range("I17").Value = 307.6
$ws.Range("J17").Value = 383.16666
$ws.Range("K17").Value = 922.8000000000001
$ws.Range("L17").Value = 1149.49998
$ws.Range("M17").Value = -753.8000000000001
$ws.Range("N17").Value = -1487.49998
$ws.Range("H122").Value = 1151.1111
$ws.Range("J122").Value = 2250
$ws.Range("L122").Value = 20250
$ws.Range("N122").Value = -25150
$ws.Range("H131").Value = 1646.8182
$ws.Range("J131").Value = 1766.5714
$ws.Range("L131").Value = 5299.7142
$ws.Range("N131").Value = -15379.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 16131152
$ws.Range("I122").Value = 1920.3914
$ws.Range("K122").Value = 5761.174199999999
$ws.Range("M122").Value = -3311.174199999999
$ws.Range("H132").Value = 4476.561
$ws.Range("I132").Value = 4426.8286
$ws.Range("J132").Value = 4766.6665
$ws.Range("K132").Value = 13280.4858
$ws.Range("L132").Value = 14299.9995
$ws.Range("M132").Value = -10750.4858
$ws.Range("N132").Value = -19359.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 10000.667
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 10000.667
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 10000.667
$ws.Range("N4").Value = -10226.667
$ws.Range("H7").Value = 4666.1665
$ws.Range("I7").Value = 4452
$ws.Range("K7").Value = 4452
$ws.Range("M7").Value = -4340
$ws.Range("H22").Value = 2047.1428
$ws.Range("I22").Value = 1807.8
$ws.Range("J22").Value = 2645.5
$ws.Range("K22").Value = 1807.8
$ws.Range("L22").Value = 2645.5
$ws.Range("M22").Value = -1512.8
$ws.Range("N22").Value = -3235.5
$ws.Range("H27").Value = 2047.1428
$ws.Range("I27").Value = 1807.8
$ws.Range("J27").Value = 2645.5
$ws.Range("K27").Value = 1807.8
$ws.Range("L27").Value = 2645.5
$ws.Range("M27").Value = -1700.8
$ws.Range("N27").Value = -2859.5
$ws.Range("H28").Value = 10000.667
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 10000.667
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 10000.667
$ws.Range("N28").Value = -10464.667
$ws.Range("H37").Value = 10000.667
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 10000.667
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 10000.667
$ws.Range("N37").Value = -10214.667
$ws.Range("H40").Value = 6457.2
$ws.Range("I40").Value = 5571.75
$ws.Range("K40").Value = 5571.75
$ws.Range("M40").Value = -5435.75
$ws.Range("H61").Value = 10958.8
$ws.Range("I61").Value = 9400
$ws.Range("K61").Value = 9400
$ws.Range("M61").Value = -9198
$ws.Range("H82").Value = 2475.0833
$ws.Range("I82").Value = 2485.4285
$ws.Range("J82").Value = 2460.6
$ws.Range("K82").Value = 2485.4285
$ws.Range("L82").Value = 2460.6
$ws.Range("M82").Value = -2124.4285
$ws.Range("N82").Value = -3182.6
$ws.Range("H85").Value = 2475.0833
$ws.Range("I85").Value = 2485.4285
$ws.Range("J85").Value = 2460.6
$ws.Range("K85").Value = 2485.4285
$ws.Range("L85").Value = 2460.6
$ws.Range("M85").Value = -1237.4285
$ws.Range("N85").Value = -4956.6
$ws.Range("H100").Value = 5498.476
$ws.Range("I100").Value = 3819.5
$ws.Range("K100").Value = 3819.5
$ws.Range("M100").Value = -3278.5
$ws.Range("H113").Value = 10958.8
$ws.Range("I113").Value = 9400
$ws.Range("K113").Value = 9400
$ws.Range("M113").Value = -7230
$ws.Range("H122").Value = 3293369.5
$ws.Range("J122").Value = 12504409
$ws.Range("L122").Value = 37513227
$ws.Range("N122").Value = -37518127
$ws.Range("H126").Value = 4666.1665
$ws.Range("I126").Value = 4452
$ws.Range("K126").Value = 13356
$ws.Range("M126").Value = -10886
$ws.Range("M4").ClearContents()
$ws.Range("M28").ClearContents()
$ws.Range("M37").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 48096596
$ws.Range("I100").Value = 72144424
$ws.Range("J100").Value = 950.5714
$ws.Range("K100").Value = 144288848
$ws.Range("L100").Value = 1901.1428
$ws.Range("M100").Value = -144288307
$ws.Range("N100").Value = -2983.1428
$ws.Range("H113").Value = 1130.5555
$ws.Range("I113").Value = 1084.375
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 3253.125
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -1083.125
$ws.Range("N113").Value = -8840
